$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.771.02"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "2.081.45"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "234.65"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -0.60%  "

$ws.Range("D7").Value = "58.73"
$ws.Range("E7").Value = "  +2.39%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  -1.64%  "

$ws.Range("D10").Value = "0.0791"
$ws.Range("E10").Value = "  +2.28%  "

$ws.Range("E11").Value = "  +2.83%  "

$ws.Range("D12").Value = "2.386.99"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").Value = "14.59"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").Value = "21.22"
$ws.Range("E14").Value = "  +3.33%  "

$ws.Range("D15").Value = "0.768"
$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").Value = "5.31"
$ws.Range("E16").Value = "  +2.62%  "

$ws.Range("D17").Value = "2.080.43"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "37.684.67"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").Value = "6.21"
$ws.Range("E19").Value = "  -3.03%  "

$ws.Range("D20").Value = "71.55"
$ws.Range("E20").Value = "  +2.49%  "

$ws.Range("D21").Value = "0.0₃0829"
$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("D22").Value = "228.47"
$ws.Range("E22").Value = "  +0.61%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  +0.41%  "

$ws.Range("D26").Value = "170.19"
$ws.Range("E26").Value = "  +2.13%  "

$ws.Range("E27").Value = "  +9.24%  "

$ws.Range("D28").Value = "9.02"
$ws.Range("E28").Value = "  +2.01%  "

$ws.Range("D29").Value = "1.43"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").Value = "19.45"
$ws.Range("E30").Value = "  +1.66%  "

$ws.Range("E31").Value = "  +1.58%  "

$ws.Range("D32").Value = "4.71"
$ws.Range("E32").Value = "  +3.91%  "

$ws.Range("D33").Value = "0.0630"
$ws.Range("E33").Value = "  +1.58%  "

$ws.Range("D34").Value = "4.68"
$ws.Range("E34").Value = "  +2.06%  "

$ws.Range("D35").Value = "2.51"
$ws.Range("E35").Value = "  +1.03%  "

$ws.Range("E36").Value = "  +6.66%  "

$ws.Range("D37").Value = "1.83"
$ws.Range("E37").Value = "  +2.52%  "

$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").Value = "5.43"
$ws.Range("E39").Value = "  -5.16%  "

$ws.Range("D40").Value = "0.0987"
$ws.Range("E40").Value = "  +1.86%  "

$ws.Range("D41").Value = "98.89"
$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("E43").Value = "  +0.75%  "

$ws.Range("D44").Value = "1.463.28"
$ws.Range("E44").Value = "  -1.53%  "

$ws.Range("D45").Value = "1.17"
$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("D46").Value = "4.25"
$ws.Range("E46").Value = "  +4.57%  "

$ws.Range("D47").Value = "16.04"
$ws.Range("E47").Value = "  +4.41%  "

$ws.Range("D48").Value = "1.06"
$ws.Range("E48").Value = "  +3.10%  "

$ws.Range("D49").Value = "7.45"
$ws.Range("E49").Value = "  +2.64%  "

$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("D51").Value = "47.53"
$ws.Range("E51").Value = "  +6.31%  "
